{"js": "// Update the date line and all twenty-five \"three-digit \u00f7 one-digit\"\n// practice answers in the table to their new values.\nconst replacements = [\n  [\"2026-01-23 Friday\", \"2026-01-24 Saturday\"],\n  [\"284\u00f77=40, 4\", \"422\u00f79=46, 8\"],\n  [\"177\u00f77=25, 2\", \"355\u00f76=59, 1\"],\n  [\"245\u00f77=35, 0\", \"206\u00f72=103, 0\"],\n  [\"577\u00f72=288, 1\", \"403\u00f78=50, 3\"],\n  [\"389\u00f73=129, 2\", \"915\u00f72=457, 1\"],\n  [\"699\u00f78=87, 3\", \"311\u00f75=62, 1\"],\n  [\"894\u00f73=298, 0\", \"695\u00f72=347, 1\"],\n  [\"444\u00f79=49, 3\", \"622\u00f78=77, 6\"],\n  [\"666\u00f76=111, 0\", \"766\u00f77=109, 3\"],\n  [\"182\u00f79=20, 2\", \"612\u00f77=87, 3\"],\n  [\"432\u00f76=72, 0\", \"245\u00f72=122, 1\"],\n  [\"691\u00f76=115, 1\", \"165\u00f73=55, 0\"],\n  [\"108\u00f76=18, 0\", \"140\u00f78=17, 4\"],\n  [\"595\u00f72=297, 1\", \"745\u00f74=186, 1\"],\n  [\"878\u00f79=97, 5\", \"896\u00f74=224, 0\"],\n  [\"291\u00f74=72, 3\", \"645\u00f78=80, 5\"],\n  [\"385\u00f75=77, 0\", \"203\u00f74=50, 3\"],\n  [\"548\u00f75=109, 3\", \"191\u00f75=38, 1\"],\n  [\"378\u00f78=47, 2\", \"896\u00f76=149, 2\"],\n  [\"172\u00f79=19, 1\", \"143\u00f74=35, 3\"],\n  [\"770\u00f77=110, 0\", \"691\u00f76=115, 1\"],\n  [\"805\u00f73=268, 1\", \"165\u00f74=41, 1\"],\n  [\"165\u00f72=82, 1\", \"983\u00f77=140, 3\"],\n  [\"743\u00f77=106, 1\", \"738\u00f73=246, 0\"],\n  [\"245\u00f72=122, 1\", \"222\u00f78=27, 6\"],\n];\n\nconst body = context.document.body;\n\n// Every \"old\" value is unique in the document, but some \"new\" values\n// happen to equal a DIFFERENT entry's \"old\" value later in the list\n// (e.g. the new text for one cell matches the old text of another).\n// To avoid re-matching an already-updated cell, resolve every search\n// BEFORE performing any of the replacements.\nconst ranges = [];\nfor (const [oldText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  ranges.push(results);\n}\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [oldText, newText] = replacements[i];\n  const items = ranges[i].items;\n  if (items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  for (const range of items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\nawait context.sync();\n", "ps1": "# Update the date line and all twenty-five \"three-digit \u00f7 one-digit\"\n# practice answers in the table to their new values.\n#\n# Cells are addressed by (row, column) position rather than by searching\n# for the old text, because a handful of the new values collide with\n# other cells' old values (e.g. the new text written into one cell is\n# the same as the old text still waiting to be replaced in a later\n# cell) \u2014 a sequential Find/Replace-by-text pass would risk rewriting\n# the wrong occurrence.\n\n$d = $word.ActiveDocument\n\n# Title line above the table.\n$d.Paragraphs.Item(1).Range.Text = \"2026-01-24 Saturday\"\n\n$table = $d.Tables.Item(1)\n\n# Each inner array is one Word table row (1-based), holding the five\n# (old, new) pairs left-to-right for that row, purely for readability /\n# sanity-checking; only the \"new\" values are actually used below.\n$rowsData = @(\n    @{ Row = 1;  Cells = @(\"422\u00f79=46, 8\", \"355\u00f76=59, 1\", \"206\u00f72=103, 0\", \"403\u00f78=50, 3\", \"915\u00f72=457, 1\") },\n    @{ Row = 5;  Cells = @(\"311\u00f75=62, 1\", \"695\u00f72=347, 1\", \"622\u00f78=77, 6\", \"766\u00f77=109, 3\", \"612\u00f77=87, 3\") },\n    @{ Row = 9;  Cells = @(\"245\u00f72=122, 1\", \"165\u00f73=55, 0\", \"140\u00f78=17, 4\", \"745\u00f74=186, 1\", \"896\u00f74=224, 0\") },\n    @{ Row = 13; Cells = @(\"645\u00f78=80, 5\", \"203\u00f74=50, 3\", \"191\u00f75=38, 1\", \"896\u00f76=149, 2\", \"143\u00f74=35, 3\") },\n    @{ Row = 17; Cells = @(\"691\u00f76=115, 1\", \"165\u00f74=41, 1\", \"983\u00f77=140, 3\", \"738\u00f73=246, 0\", \"222\u00f78=27, 6\") }\n)\n\nforeach ($rowInfo in $rowsData) {\n    $r = $rowInfo.Row\n    for ($c = 1; $c -le 5; $c++) {\n        $table.Cell($r, $c).Range.Text = $rowInfo.Cells[$c - 1]\n    }\n}\n"}
